# Insert a new data row above the current row 45 (old row 45 and everything
# below it shifts down by one, growing the used range from A1:R76 to A1:R77),
# then populate the newly-inserted row 45 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 45..76 down to 46..77, leaving a blank row 45.
$ws.Rows.Item(45).Insert()

# Fill the new row 45.
$ws.Cells.Item(45, 1).Value  = 4
$ws.Cells.Item(45, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(45, 3).Value  = "Los Lagos"
$ws.Cells.Item(45, 4).Value  = 44574
$ws.Cells.Item(45, 5).Value  = 10
$ws.Cells.Item(45, 6).Value  = 100112026
$ws.Cells.Item(45, 7).Value  = "Haba"
$ws.Cells.Item(45, 8).Value  = "Sin especificar"
$ws.Cells.Item(45, 9).Value  = "Primera"
$ws.Cells.Item(45, 10).Value = 70
$ws.Cells.Item(45, 11).Value = 24000
$ws.Cells.Item(45, 12).Value = 24000
$ws.Cells.Item(45, 13).Value = 24000
$ws.Cells.Item(45, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(45, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(45, 16).Value = 960
$ws.Cells.Item(45, 17).Value = 25
$ws.Cells.Item(45, 18).Value = "Hortaliza"
